$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write column A (Sending cluster) for all data rows first to control shared-string order
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "FAPs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("A6").Value = "MuSCs"
$ws.Range("A7").Value = "MuSCs"

# Write column B (Ligand symbol) for all data rows
$ws.Range("B2").Value = "Gm13306"
$ws.Range("B3").Value = "Gm13306"
$ws.Range("B4").Value = "Gm13306"
$ws.Range("B5").Value = "Gm13306"
$ws.Range("B6").Value = "Gm13306"
$ws.Range("B7").Value = "Gm13306"

# Write column C (Receptor symbol) for all data rows
$ws.Range("C2").Value = "Ccr10"
$ws.Range("C3").Value = "Ccr10"
$ws.Range("C4").Value = "Ccr10"
$ws.Range("C5").Value = "Ccr10"
$ws.Range("C6").Value = "Ccr10"
$ws.Range("C7").Value = "Ccr10"

# Write column D (Target cluster) for all data rows
$ws.Range("D2").Value = "FAPs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "MuSCs"

# Write numeric columns E through T for all data rows
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.065096
$ws.Range("H2").Value = 0.195288
$ws.Range("I2").Value = 0.02794828919627058
$ws.Range("J2").Value = 0.02794828919627058
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.036595333333333
$ws.Range("N2").Value = 3.109786
$ws.Range("O2").Value = 0.393072250513715
$ws.Range("P2").Value = 0.3930722505137151
$ws.Range("Q2").Value = 0.06747820981866666
$ws.Range("R2").Value = 0.607303888368
$ws.Range("S2").Value = 0.01098569693238623
$ws.Range("T2").Value = 0.01098569693238623
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.065096
$ws.Range("H3").Value = 0.195288
$ws.Range("I3").Value = 0.02794828919627058
$ws.Range("J3").Value = 0.02794828919627058
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.600567
$ws.Range("N3").Value = 4.801701
$ws.Range("O3").Value = 0.6069277494862849
$ws.Range("P3").Value = 0.6069277494862849
$ws.Range("Q3").Value = 0.104190509432
$ws.Range("R3").Value = 0.9377145848879999
$ws.Range("S3").Value = 0.01696259226388436
$ws.Range("T3").Value = 0.01696259226388436
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.553094
$ws.Range("H4").Value = 4.659282
$ws.Range("I4").Value = 0.6668047231933247
$ws.Range("J4").Value = 0.6668047231933247
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.036595333333333
$ws.Range("N4").Value = 3.109786
$ws.Range("O4").Value = 0.393072250513715
$ws.Range("P4").Value = 0.3930722505137151
$ws.Range("Q4").Value = 1.609929992628
$ws.Range("R4").Value = 14.489369933652
$ws.Range("S4").Value = 0.262102433198775
$ws.Range("T4").Value = 0.262102433198775
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.553094
$ws.Range("H5").Value = 4.659282
$ws.Range("I5").Value = 0.6668047231933247
$ws.Range("J5").Value = 0.6668047231933247
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.600567
$ws.Range("N5").Value = 4.801701
$ws.Range("O5").Value = 0.6069277494862849
$ws.Range("P5").Value = 0.6069277494862849
$ws.Range("Q5").Value = 2.485831004298
$ws.Range("R5").Value = 22.372479038682
$ws.Range("S5").Value = 0.4047022899945498
$ws.Range("T5").Value = 0.4047022899945498
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7109686666666667
$ws.Range("H6").Value = 2.132906
$ws.Range("I6").Value = 0.3052469876104046
$ws.Range("J6").Value = 0.3052469876104047
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.036595333333333
$ws.Range("N6").Value = 3.109786
$ws.Range("O6").Value = 0.393072250513715
$ws.Range("P6").Value = 0.3930722505137151
$ws.Range("Q6").Value = 0.7369868020128889
$ws.Range("R6").Value = 6.632881218116001
$ws.Range("S6").Value = 0.1199841203825538
$ws.Range("T6").Value = 0.1199841203825539
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7109686666666667
$ws.Range("H7").Value = 2.132906
$ws.Range("I7").Value = 0.3052469876104046
$ws.Range("J7").Value = 0.3052469876104047
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.600567
$ws.Range("N7").Value = 4.801701
$ws.Range("O7").Value = 0.6069277494862849
$ws.Range("P7").Value = 0.6069277494862849
$ws.Range("Q7").Value = 1.137952985900667
$ws.Range("R7").Value = 10.241576873106
$ws.Range("S7").Value = 0.1852628672278508
$ws.Range("T7").Value = 0.1852628672278508
